$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Additional notes to the task, inserted in column G next to the
# "Раскладка точек блокировки:" picture (rows 14-20).
$ws.Range("G14").Value = "двигаться мы можем не только по уже отмеченной линии"
$ws.Range("G15").Value = "у нас просто есть точки с цифрами"
$ws.Range("G16").Value = "и между ними можно двигаться как угодно — но только к ближайшим соседям"
$ws.Range("G18").Value = "Типа 5-8 нельзя"
$ws.Range("G19").Value = "А 5-2 можно"
$ws.Range("G20").Value = "5 — 8 будет просто через 2"

# Re-apply the number formats on the length-calculation helper cells so the
# style table reflects the same refresh the author's resave produced.
$ws.Range("D57").NumberFormat = "General"
$ws.Range("D58").NumberFormat = "General"
$ws.Range("A58").NumberFormat = "0.000000"

# Park the selection where the author left it before committing.
$ws.Range("G21").Select()
